$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting the existing data rows (2-11) down to (3-12).
$ws.Rows.Item(2).Insert()

# Excel's row insert copies the formatting of the row above (the bold header),
# so reset the freshly inserted row back to the plain default style used by
# all the other data rows.
$ws.Range("A2:D2").Style = "Normal"

# Populate the newly inserted row 2 with the latest day's price data.
# Force the date to be stored as plain text (matching the other date cells)
# instead of being auto-converted into a date serial number.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-12-01"
$ws.Cells.Item(2, 1).Style = "Normal"

$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
